$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 ("2021年") with the same formatting as the previous year
# rows (A2:A4) by copying the formatting from A4, then filling in the text.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").Value = "2021年"

$ws.Range("B5").Value = 2.9
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 6.2
